$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCalc  = $wb.Worksheets.Item("Calculations")
$wsEhp   = $wb.Worksheets.Item("EHPpUC")

# ---------------------------------------------------------------------------
# Bug fix: the "Hydrogen Energy Density" input cell (lb per metric ton of H2)
# had the wrong value. Insert a blank row below it (room for the note block
# further down) and correct the figure.
# ---------------------------------------------------------------------------
$wsCalc.Rows.Item(9).Insert()
$wsCalc.Range("A8").Value = 60920

# Explanation / notes about where the corrected number comes from, placed to
# the right of the "lb per metric ton" / "BTU H2/MW" block.
$wsCalc.Range("C12").Value = "The study they are citing uses 3 scenarios of Fuel Cell Electric Vehicle adoption. Then it calculates the amount of hydrogen needed to suppor those vehicles. Then it calculates the electrolyzer capacity needed to supply that hydrogen. "
$wsCalc.Range("C13").Value = "So, I think it's fair, using the EPS assumptions of 24/7/365 operation, that smallest electrolyzer you would need to produce 1.39e10 annual Btu would be 1 MW."
$wsCalc.Range("C14").Value = "No reason to think this would be different for Texas."

$wsCalc.Range("C12:C14").Font.Color = 12874308

# Restore selections roughly where they were left when the file was saved.
$wsCalc.Range("J6").Select()
$wsEhp.Range("B2").Select()
$wsAbout.Range("B6").Select()
